$wb = $excel.ActiveWorkbook

# Overview sheet: bump the "Latest HO Xliff Generate Date" for the
# 3c479a4d... row (handback report regenerated).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-02 14:59:51"

# zh-cn sheet: refresh handoff/handback datetimes and flag metadata present.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-02 14:59:47"
$zhcn.Range("K2").Value = "2016-09-02 15:00:40"
$zhcn.Range("O3").Value = "'True"

# de-de sheet: refresh handoff/handback datetimes.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-02 14:59:51"
$dede.Range("K2").Value = "2016-09-02 15:00:57"
